# Updates the cryptocurrency price/volume snapshot on Sheet1 to the
# latest scrape, per the GitHub Actions cron commit on
# Thu Sep  7 07:40:13 UTC 2023.
#
# Columns: D = Price (text), E = Volume(1h) (text, e.g. "  +0.34%  ").
# A handful of new Price values are plain decimals (e.g. "215.69")
# that Excel would otherwise auto-convert to a Number; those are
# written with a leading apostrophe so COM keeps them as Text,
# matching the original inline-string cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.788.17"

$ws.Range("D3").Value = "1.638.12"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'215.69"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("E6").Value = "  -0.61%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D10").Value = "'19.64"
$ws.Range("E10").Value = "  -1.34%  "

$ws.Range("D11").Value = "'0.0794"
$ws.Range("E11").Value = "  +1.67%  "

$ws.Range("D12").Value = "'4.27"
$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("D13").Value = "1.864.16"
$ws.Range("E13").Value = "  +0.31%  "

$ws.Range("D14").Value = "1.636.41"
$ws.Range("E14").Value = "  -0.22%  "

$ws.Range("E15").Value = "  +0.96%  "

$ws.Range("D16").Value = "0.0₃0764"
$ws.Range("E16").Value = "  -0.57%  "

$ws.Range("D17").Value = "'63.21"
$ws.Range("E17").Value = "  +0.28%  "

$ws.Range("D18").Value = "25.821.53"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("E20").Value = "  +2.14%  "

$ws.Range("D21").Value = "'192.61"
$ws.Range("E21").Value = "  -0.89%  "

$ws.Range("E22").Value = "  +0.48%  "

$ws.Range("D23").Value = "'6.30"
$ws.Range("E23").Value = "  +1.70%  "

$ws.Range("D24").Value = "'1.84"
$ws.Range("E24").Value = "  +4.65%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").Value = "'141.70"
$ws.Range("E26").Value = "  +1.43%  "

$ws.Range("E27").Value = "  +1.21%  "

$ws.Range("E28").Value = "  +1.44%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("E31").Value = "  -0.48%  "

$ws.Range("E32").Value = "  +0.65%  "

$ws.Range("E33").Value = "  -0.64%  "

$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("E36").Value = "  +0.39%  "

$ws.Range("D37").Value = "1.136.22"
$ws.Range("E37").Value = "  +1.23%  "

$ws.Range("E38").Value = "  -1.44%  "

$ws.Range("E39").Value = "  -0.88%  "

$ws.Range("E40").Value = "  -0.22%  "

$ws.Range("E41").Value = "  +0.19%  "

$ws.Range("D42").Value = "'5.58"
$ws.Range("E42").Value = "  +0.82%  "

$ws.Range("D43").Value = "'100.67"
$ws.Range("E43").Value = "  +0.99%  "

$ws.Range("D44").Value = "'0.802"
$ws.Range("E44").Value = "  +0.32%  "

$ws.Range("D45").Value = "1.773.22"
$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("E46").Value = "  +3.44%  "

$ws.Range("D47").Value = "'55.35"

$ws.Range("E49").Value = "  -0.16%  "

$ws.Range("E50").Value = "  +3.76%  "

$ws.Range("E51").Value = "  -1.79%  "
